$d = $word.ActiveDocument

# 1. Section 7.1 heading: "Thận trọng trong việc bảo quản" -> "Thận trọng trong xử lý an toàn"
$d.Content.Find.Execute("Thận trọng trong việc bảo quản", $true, $false, $false, $false, $false, $true, 1, $false, "Thận trọng trong xử lý an toàn", 2)

# 2. "Hơi tương đối" -> "Tỉ trọng hơi"
$d.Content.Find.Execute("Hơi tương đối", $true, $false, $false, $false, $false, $true, 1, $false, "Tỉ trọng hơi", 2)

# 3. "Tính chất oxy hóa" -> "Tính oxy hóa"
$d.Content.Find.Execute("Tính chất oxy hóa", $true, $false, $false, $false, $false, $true, 1, $false, "Tính oxy hóa", 2)

# 4. "bon dioxide" -> "bon đioxit" (Cacbon dioxide -> Cacbon đioxit)
$d.Content.Find.Execute("bon dioxide", $true, $false, $false, $false, $false, $true, 1, $false, "bon đioxit", 2)

# 5. "bon monoxide" -> "bon monoxit" (cacbon monoxide -> cacbon monoxit)
$d.Content.Find.Execute("bon monoxide", $true, $false, $false, $false, $false, $true, 1, $false, "bon monoxit", 2)

# 6. "Độc tính cấp (đường miệng)" -> "Độc cấp tính (đường miệng)"
$d.Content.Find.Execute("Độc tính cấp (đường miệng)", $true, $false, $false, $false, $false, $true, 1, $false, "Độc cấp tính (đường miệng)", 2)

# 7. "Độc tính cấp (qua da)" -> "Độc cấp tính (qua da)"
$d.Content.Find.Execute("Độc tính cấp (qua da)", $true, $false, $false, $false, $false, $true, 1, $false, "Độc cấp tính (qua da)", 2)
